$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row (Joel Juaristi) first, so the shared-string table
# fills in "Joel" / "Juaristi" / the e-mail before the header text is
# rewritten (keeps the shared-string order identical to a real Excel save).

# Numero de socio (numeric id column) - reuse formatting from A2 (s=2)
$ws.Range("A7").Value = 6
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# Nombre
$ws.Range("B7").Value = "Joel"
$ws.Range("B2").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Apellidos
$ws.Range("C7").Value = "Juaristi"
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)

# Correo - becomes a mailto hyperlink, still displaying the e-mail address
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:joeljuaristi@hotmail.com")
$ws.Range("D7").Value = "joeljuaristi@hotmail.com"

# Match the explicit row height used for the new row
$ws.Rows.Item(7).RowHeight = 15.75

# --- Re-label the header row (now lower-case) and move the "numero_socio"
# header to column A, matching the re-ordered headers.
$ws.Range("B1").Value = "nombre"
$ws.Range("C1").Value = "apellidos"
$ws.Range("D1").Value = "correo"
$ws.Range("A1").Value = "numero_socio"

Write-Host "Workbook updated"
